$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 165 currently has empty E165/F165 inline-string cells; clear them out.
$ws.Range("E165").Value = $null
$ws.Range("F165").Value = $null

# Append the new attendance records collected in the "Diario de Classe" screen.
$newRows = @(
    @(2, "maria", "2024-03-01", "A", $null, $null),
    @(2, "maria", "2024-01-10", "P", "Não tem", "Também não tem"),
    @(2, "maria", "2024-10-07", "P", "Testeeeee", "Testeeee"),
    @(2, "maria", "2024-10-07", "P", "Hoje teve aula", "Hoje Teve de fato aula")
)

$r = 166
foreach ($row in $newRows) {
    # Every column in this sheet stores plain text (inline strings), even the
    # numeric-looking id column and the yyyy-mm-dd date column, so force text
    # number formatting before writing to stop Excel auto-converting them to
    # a real number / date.
    $rowRange = $ws.Range("A$r`:F$r")
    $rowRange.NumberFormat = "@"

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    if ($row[4] -ne $null) { $ws.Cells.Item($r, 5).Value = $row[4] }
    if ($row[5] -ne $null) { $ws.Cells.Item($r, 6).Value = $row[5] }
    $r++
}
